$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data preprocessing results (recomputed values) in column C
$ws.Range("C3").Value = -172.354577427518
$ws.Range("C7").Value = 75.625968424842398
$ws.Range("C9").Value = 66.5107695914093
$ws.Range("C10").Value = 88.329189932518304
$ws.Range("C11").Value = 84.016662136192295
$ws.Range("C12").Value = 86.770905935039295

# Move active selection to D9, matching the saved cursor position
$ws.Activate()
$ws.Range("D9").Select()
